$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" column header in F1, matching the style of the other headers
# (copy E1's formatting -- bold/centered/bordered header style -- onto F1, then set its text)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamps for each data row (F2:F19)
$timestamps = @(
    "2021-10-05 10:52:57.379349",
    "2021-10-05 10:52:57.379361",
    "2021-10-05 10:52:57.379365",
    "2021-10-05 10:52:57.379368",
    "2021-10-05 10:52:57.379372",
    "2021-10-05 10:52:57.379375",
    "2021-10-05 10:52:57.379378",
    "2021-10-05 10:52:57.379381",
    "2021-10-05 10:52:57.379384",
    "2021-10-05 10:52:57.379388",
    "2021-10-05 10:52:57.379391",
    "2021-10-05 10:52:57.379394",
    "2021-10-05 10:52:57.379397",
    "2021-10-05 10:52:57.379400",
    "2021-10-05 10:52:57.379403",
    "2021-10-05 10:52:57.379406",
    "2021-10-05 10:52:57.379409",
    "2021-10-05 10:52:57.379412"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
